$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.050050253866516
$ws.Cells.Item(2, 3).Value = 0.2310093480059265
$ws.Cells.Item(2, 5).Value = 0.5929007669134876
$ws.Cells.Item(2, 6).Value = 2.3791024842524
$ws.Cells.Item(2, 7).Value = 0.6208773747725473
$ws.Cells.Item(2, 8).Value = 0.6856196669878329
$ws.Cells.Item(2, 10).Value = 0.05044172146965664
$ws.Cells.Item(2, 14).Value = 0.9441726390137362

$ws.Cells.Item(3, 2).Value = 0.9352176752354353
$ws.Cells.Item(3, 3).Value = 0.2012804013143352
$ws.Cells.Item(3, 5).Value = 0.5768069957385507
$ws.Cells.Item(3, 6).Value = 2.331772135540007
$ws.Cells.Item(3, 7).Value = 0.6064439058979758
$ws.Cells.Item(3, 8).Value = 0.6850004223734345
$ws.Cells.Item(3, 10).Value = 0.04984953454162877
$ws.Cells.Item(3, 14).Value = 0.9573717346602706

$ws.Cells.Item(4, 2).Value = 0.8648972588064794
$ws.Cells.Item(4, 3).Value = 0.1829853590511448
$ws.Cells.Item(4, 5).Value = 0.5671837325725022
$ws.Cells.Item(4, 6).Value = 2.304323275465237
$ws.Cells.Item(4, 7).Value = 0.5982340069400038
$ws.Cells.Item(4, 8).Value = 0.6851366307958244
$ws.Cells.Item(4, 10).Value = 0.04953422531249174
$ws.Cells.Item(4, 14).Value = 0.9659772779770215

$ws.Cells.Item(5, 2).Value = 0.8362881110776357
$ws.Cells.Item(5, 3).Value = 0.1755195289977962
$ws.Cells.Item(5, 5).Value = 0.5633270900993068
$ws.Cells.Item(5, 6).Value = 2.293541300788732
$ws.Cells.Item(5, 7).Value = 0.5950511024394558
$ws.Cells.Item(5, 8).Value = 0.6853215024652144
$ws.Cells.Item(5, 10).Value = 0.04941779300844829
$ws.Cells.Item(5, 14).Value = 0.9696100427429712

$ws.Cells.Item(6, 2).Value = 0.8315404199576051
$ws.Cells.Item(6, 3).Value = 0.1742792008863887
$ws.Cells.Item(6, 5).Value = 0.5626906167700909
$ws.Cells.Item(6, 6).Value = 2.29177529179772
$ws.Cells.Item(6, 7).Value = 0.5945323711357844
$ws.Cells.Item(6, 8).Value = 0.6853599992424648
$ws.Cells.Item(6, 10).Value = 0.04939918533725418
$ws.Cells.Item(6, 14).Value = 0.9702208638159604

$ws.Cells.Item(7, 2).Value = 0.8645112356727509
$ws.Cells.Item(7, 3).Value = 0.1828847146394708
$ws.Cells.Item(7, 5).Value = 0.5671314577738116
$ws.Cells.Item(7, 6).Value = 2.304176233884718
$ws.Cells.Item(7, 7).Value = 0.5981904241921541
$ws.Cells.Item(7, 8).Value = 0.6851386009035707
$ws.Cells.Item(7, 10).Value = 0.0495326063522441
$ws.Cells.Item(7, 14).Value = 0.9660257609473852

$ws.Cells.Item(8, 2).Value = 1.010417019938927
$ws.Cells.Item(8, 3).Value = 0.2207673463941546
$ws.Cells.Item(8, 5).Value = 0.5872979674275669
$ws.Cells.Item(8, 6).Value = 2.362447306122348
$ws.Cells.Item(8, 7).Value = 0.6157645184207183
$ws.Cells.Item(8, 8).Value = 0.6852986662868403
$ws.Cells.Item(8, 10).Value = 0.05022745959755071
$ws.Cells.Item(8, 14).Value = 0.9486195889842968

$ws.Cells.Item(9, 2).Value = 1.298043484719358
$ws.Cells.Item(9, 3).Value = 0.2947332587150413
$ws.Cells.Item(9, 5).Value = 0.6288998994232031
$ws.Cells.Item(9, 6).Value = 2.489596343397437
$ws.Cells.Item(9, 7).Value = 0.6554659106669476
$ws.Cells.Item(9, 8).Value = 0.6897336610103935
$ws.Cells.Item(9, 10).Value = 0.05197729242039628
$ws.Cells.Item(9, 14).Value = 0.9184680104150615

$ws.Cells.Item(10, 2).Value = 1.51033079686448
$ws.Cells.Item(10, 3).Value = 0.3488940257901163
$ws.Cells.Item(10, 5).Value = 0.6607294495293985
$ws.Cells.Item(10, 6).Value = 2.591000059584331
$ws.Cells.Item(10, 7).Value = 0.687919151624385
$ws.Cells.Item(10, 8).Value = 0.695538642732771
$ws.Cells.Item(10, 10).Value = 0.05350475428127766
$ws.Cells.Item(10, 14).Value = 0.8987495243295172

$ws.Cells.Item(11, 2).Value = 1.607128475660147
$ws.Cells.Item(11, 3).Value = 0.373497100714701
$ws.Cells.Item(11, 5).Value = 0.6754869515463184
$ws.Cells.Item(11, 6).Value = 2.638895997425436
$ws.Cells.Item(11, 7).Value = 0.7034161458230415
$ws.Cells.Item(11, 8).Value = 0.6987400097721093
$ws.Cells.Item(11, 10).Value = 0.05425342547339795
$ws.Cells.Item(11, 14).Value = 0.8903090979956474

$ws.Cells.Item(12, 2).Value = 1.643816446109952
$ws.Cells.Item(12, 3).Value = 0.3828087895452086
$ws.Cells.Item(12, 5).Value = 0.681115374457363
$ws.Cells.Item(12, 6).Value = 2.65728924739318
$ws.Cells.Item(12, 7).Value = 0.7093915334055509
$ws.Cells.Item(12, 8).Value = 0.7000334988406962
$ws.Cells.Item(12, 10).Value = 0.05454476730794511
$ws.Cells.Item(12, 14).Value = 0.8871892462667503

$ws.Cells.Item(13, 2).Value = 1.635913574579035
$ws.Cells.Item(13, 3).Value = 0.3808035703910377
$ws.Cells.Item(13, 5).Value = 0.679901408948254
$ws.Cells.Item(13, 6).Value = 2.653316509398763
$ws.Cells.Item(13, 7).Value = 0.7080998436225343
$ws.Cells.Item(13, 8).Value = 0.6997513021019302
$ws.Cells.Item(13, 10).Value = 0.05448167155808648
$ws.Cells.Item(13, 14).Value = 0.8878577630412963

$ws.Cells.Item(14, 2).Value = 1.610146159079591
$ws.Cells.Item(14, 3).Value = 0.3742632780289341
$ws.Cells.Item(14, 5).Value = 0.6759492015469419
$ws.Cells.Item(14, 6).Value = 2.640404076908595
$ws.Cells.Item(14, 7).Value = 0.7039055913293737
$ws.Cells.Item(14, 8).Value = 0.6988447951694354
$ws.Cells.Item(14, 10).Value = 0.05427723673250995
$ws.Cells.Item(14, 14).Value = 0.8900508944998435

$ws.Cells.Item(15, 2).Value = 1.594367147568562
$ws.Cells.Item(15, 3).Value = 0.3702565158585571
$ws.Cells.Item(15, 5).Value = 0.6735335823830866
$ws.Cells.Item(15, 6).Value = 2.632528251146823
$ws.Cells.Item(15, 7).Value = 0.7013504710890857
$ws.Cells.Item(15, 8).Value = 0.6983001255493093
$ws.Cells.Item(15, 10).Value = 0.05415303790605464
$ws.Cells.Item(15, 14).Value = 0.891404200446388

$ws.Cells.Item(16, 2).Value = 1.504009421826368
$ws.Cells.Item(16, 3).Value = 0.3472854592160957
$ws.Cells.Item(16, 5).Value = 0.659770615279939
$ws.Cells.Item(16, 6).Value = 2.587905673831727
$ws.Cells.Item(16, 7).Value = 0.6869212803194671
$ws.Cells.Item(16, 8).Value = 0.6953407570504737
$ws.Cells.Item(16, 10).Value = 0.05345691814137865
$ws.Cells.Item(16, 14).Value = 0.8993117998350257

$ws.Cells.Item(17, 2).Value = 1.44863614069277
$ws.Cells.Item(17, 3).Value = 0.3331845555223936
$ws.Cells.Item(17, 5).Value = 0.651398765708052
$ws.Cells.Item(17, 6).Value = 2.560985138885911
$ws.Cells.Item(17, 7).Value = 0.6782584498738515
$ws.Cells.Item(17, 8).Value = 0.6936692874635924
$ws.Cells.Item(17, 10).Value = 0.05304372478180852
$ws.Cells.Item(17, 14).Value = 0.9042986397553747

$ws.Cells.Item(18, 2).Value = 1.41680826206067
$ws.Cells.Item(18, 3).Value = 0.3250707969374105
$ws.Cells.Item(18, 5).Value = 0.6466096677746833
$ws.Cells.Item(18, 6).Value = 2.545667343688365
$ws.Cells.Item(18, 7).Value = 0.67334480318533
$ws.Cells.Item(18, 8).Value = 0.692760638571059
$ws.Cells.Item(18, 10).Value = 0.05281112694923706
$ws.Cells.Item(18, 14).Value = 0.9072167754711344

$ws.Cells.Item(19, 2).Value = 1.40603555555515
$ws.Cells.Item(19, 3).Value = 0.3223230505607262
$ws.Cells.Item(19, 5).Value = 0.6449926533293677
$ws.Cells.Item(19, 6).Value = 2.540509482555478
$ws.Cells.Item(19, 7).Value = 0.6716929298594607
$ws.Cells.Item(19, 8).Value = 0.6924620255794025
$ws.Cells.Item(19, 10).Value = 0.05273323945075958
$ws.Cells.Item(19, 14).Value = 0.9082133598120876

$ws.Cells.Item(20, 2).Value = 1.454528508586066
$ws.Cells.Item(20, 3).Value = 0.3346859610388151
$ws.Cells.Item(20, 5).Value = 0.6522872544139631
$ws.Cells.Item(20, 6).Value = 2.563833661795911
$ws.Cells.Item(20, 7).Value = 0.6791734726411676
$ws.Cells.Item(20, 8).Value = 0.693841756094514
$ws.Cells.Item(20, 10).Value = 0.05308718552371516
$ws.Cells.Item(20, 14).Value = 0.9037626227822884

$ws.Cells.Item(21, 2).Value = 1.61771378255844
$ws.Cells.Item(21, 3).Value = 0.3761844545798567
$ws.Cells.Item(21, 5).Value = 0.6771089722386989
$ws.Cells.Item(21, 6).Value = 2.644189804658907
$ws.Cells.Item(21, 7).Value = 0.7051346287312299
$ws.Cells.Item(21, 8).Value = 0.6991088498126885
$ws.Cells.Item(21, 10).Value = 0.05433707074509897
$ws.Cells.Item(21, 14).Value = 0.8894046447761284

$ws.Cells.Item(22, 2).Value = 1.724555994607329
$ws.Cells.Item(22, 3).Value = 0.403277230741935
$ws.Cells.Item(22, 5).Value = 0.6935650458163849
$ws.Cells.Item(22, 6).Value = 2.698200880986235
$ws.Cells.Item(22, 7).Value = 0.722726190445087
$ws.Cells.Item(22, 8).Value = 0.7030247738364039
$ws.Cells.Item(22, 10).Value = 0.05519966450258096
$ws.Cells.Item(22, 14).Value = 0.8804660261602137

$ws.Cells.Item(23, 2).Value = 1.667514751658189
$ws.Cells.Item(23, 3).Value = 0.3888199279095943
$ws.Cells.Item(23, 5).Value = 0.6847607199329246
$ws.Cells.Item(23, 6).Value = 2.669236804836089
$ws.Cells.Item(23, 7).Value = 0.7132796127514496
$ws.Cells.Item(23, 8).Value = 0.70089124640009
$ws.Cells.Item(23, 10).Value = 0.05473506623414437
$ws.Cells.Item(23, 14).Value = 0.8851959366293016

$ws.Cells.Item(24, 2).Value = 1.451864547399964
$ws.Cells.Item(24, 3).Value = 0.3340071972965575
$ws.Cells.Item(24, 5).Value = 0.651885493953543
$ws.Cells.Item(24, 6).Value = 2.562545348914711
$ws.Cells.Item(24, 7).Value = 0.6787595831052897
$ws.Cells.Item(24, 8).Value = 0.6937636201539874
$ws.Cells.Item(24, 10).Value = 0.05306752150388405
$ws.Cells.Item(24, 14).Value = 0.9040047967863316

$ws.Cells.Item(25, 2).Value = 1.220066589971395
$ws.Cells.Item(25, 3).Value = 0.2747565230165776
$ws.Cells.Item(25, 5).Value = 0.6174242482057508
$ws.Cells.Item(25, 6).Value = 2.453806385166644
$ws.Cells.Item(25, 7).Value = 0.6441554256199851
$ws.Cells.Item(25, 8).Value = 0.6880891633542063
$ws.Cells.Item(25, 10).Value = 0.05146186458701862
$ws.Cells.Item(25, 14).Value = 0.9261978924210084
